$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 0.0005555555555555556
$ws.Range("K2").Value = 3939
$ws.Range("L2").Value = 0.007878
